$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 575, shifting existing rows 575-617 down to 576-618.
$ws.Rows.Item(575).Insert()

# Populate the newly inserted row 575 with its data.
$ws.Range("A575").Value = 10
$ws.Range("B575").Value = "Vega Modelo de Temuco"
$ws.Range("C575").Value = "La Araucanía"
$ws.Range("D575").Value = 44783
$ws.Range("E575").Value = 9
$ws.Range("F575").Value = 100112045
$ws.Range("G575").Value = "Zapallo"
$ws.Range("H575").Value = "Camote"
$ws.Range("I575").Value = "1a (guarda)"
$ws.Range("J575").Value = 200
$ws.Range("K575").Value = 1200
$ws.Range("L575").Value = 1200
$ws.Range("M575").Value = 1200
$ws.Range("N575").Value = "$/kilo (volumen en unidades)"
$ws.Range("O575").Value = "Región de O'Higgins"
$ws.Range("P575").Value = 1200
$ws.Range("Q575").Value = 1
$ws.Range("R575").Value = "Hortaliza"
